$d = $word.ActiveDocument

# "A planilha de riscos para incrementar"
#
# Two leftover double-spaces in the "Solução proposta" paragraph get
# tidied up to single spaces:
#   "como também  planejar o melhor"       -> "como também planejar o melhor"
#   "de espaço para o marketing,  visando" -> "de espaço para o marketing, visando"
#
# Each fix is applied as a single-character deletion of the redundant
# space (rather than a full-phrase Find/Replace) so the edit stays as
# surgical as the COM object model allows. The rightmost fix is applied
# first so it does not shift the still-to-be-computed offsets of the
# left one.

# --- Fix 2 (rightmost): "marketing,  visando" ------------------------------
$f2 = $d.Content.Duplicate
$f2.Find.Execute("marketing,  visando") | Out-Null
$gap2 = $d.Range($f2.Start + 11, $f2.Start + 12)
$gap2.Text = ""

# --- Fix 1 (leftmost): "também  planejar" -----------------------------------
$f1 = $d.Content.Duplicate
$f1.Find.Execute("também  planejar") | Out-Null
$gap1 = $d.Range($f1.Start + 7, $f1.Start + 8)
$gap1.Text = ""
